$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '42.987.63'
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").Value = '2.539.88'
$ws.Range("E3").Value = '  -0.69%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("E5").Value = '  +1.54%  '

$ws.Range("D6").Value = "'100.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.27%  '

$ws.Range("D7").Value = "'0.581"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.26%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").Value = "'0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.57%  '

$ws.Range("D10").Value = "'37.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.89%  '

$ws.Range("D11").Value = "'0.0819"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.10%  '

$ws.Range("D12").Value = "'7.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("E13").Value = '  -0.24%  '

$ws.Range("D14").Value = '2.929.77'
$ws.Range("E14").Value = '  -0.83%  '

$ws.Range("D15").Value = '2.553.22'

$ws.Range("D16").Value = "'15.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.73%  '

$ws.Range("E17").Value = '  -0.82%  '

$ws.Range("D18").Value = '43.015.09'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").Value = "'13.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.35%  '

$ws.Range("D20").Value = '0.0₃0989'
$ws.Range("E20").Value = '  -0.36%  '

$ws.Range("E21").Value = '  -0.92%  '

$ws.Range("D22").Value = "'71.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.31%  '

$ws.Range("D23").Value = "'254.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.49%  '

$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.44%  '

$ws.Range("E25").Value = '  -3.08%  '

$ws.Range("D26").Value = "'27.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.91%  '

$ws.Range("E27").Value = '  +0.23%  '

$ws.Range("D28").Value = "'10.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.40%  '

$ws.Range("E29").Value = '  +9.67%  '

$ws.Range("D30").Value = "'39.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.30%  '

$ws.Range("E31").Value = '  +2.76%  '

$ws.Range("D32").Value = "'158.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.09%  '

$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = "'2.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = "'0.0801"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("D35").Value = "'3.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.61%  '

$ws.Range("D36").Value = "'2.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.76%  '

$ws.Range("D37").Value = "'18.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.42%  '

$ws.Range("E38").Value = '  +1.97%  '

$ws.Range("E39").Value = '  +0.32%  '

$ws.Range("D40").Value = "'24.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.23%  '

$ws.Range("D41").Value = "'3.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.47%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'3.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.26%  '

$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").Value = "'2.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.89%  '

$ws.Range("E44").Value = '  -1.64%  '

$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").Value = '2.052.10'
$ws.Range("E46").Value = '  -2.05%  '

$ws.Range("D47").Value = "'86.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.24%  '

$ws.Range("D48").Value = "'9.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.05%  '

$ws.Range("D49").Value = '2.787.42'
$ws.Range("E49").Value = '  -0.82%  '

$ws.Range("E50").Value = '  +1.10%  '

$ws.Range("D51").Value = "'103.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.40%  '
